# Applies the "Updated cryptos list" data refresh described in the commit diff.
# Price (column D) and Volume(1h) (column E) values are refreshed for most rows;
# rows 32/33 additionally swap their Coin name + Link (HuobiToken <-> ImmutableX).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking text (e.g. "310.80", "0.000008713") that must
# stay plain text, matching the workbook source (inline strings, not numbers).
# Force text format while writing, then drop the format again so no stray cell
# style is introduced.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = '26.934.61'
$ws.Range("E2").Value = '  -1.38%  '
$ws.Range("D3").Value = '1.828.54'
$ws.Range("E3").Value = '  -1.52%  '
$ws.Range("E4").Value = '  +0.65%  '
$ws.Range("D5").Value = '310.80'
$ws.Range("E5").Value = '  -0.96%  '
$ws.Range("D6").Value = '1.006'
$ws.Range("E6").Value = '  +0.48%  '
$ws.Range("D7").Value = '0.4575'
$ws.Range("E7").Value = '  -0.90%  '
$ws.Range("D8").Value = '0.3704'
$ws.Range("D9").Value = '0.07170'
$ws.Range("E9").Value = '  -2.02%  '
$ws.Range("D10").Value = '0.8755'
$ws.Range("E10").Value = '  -0.79%  '
$ws.Range("D11").Value = '0.07767'
$ws.Range("E11").Value = '  -0.72%  '
$ws.Range("D12").Value = '19.62'
$ws.Range("E12").Value = '  -1.03%  '
$ws.Range("D13").Value = '1.860.15'
$ws.Range("E13").Value = '  +0.06%  '
$ws.Range("D14").Value = '5.328'
$ws.Range("E14").Value = '  -1.02%  '
$ws.Range("D15").Value = '6.394'
$ws.Range("E15").Value = '  -2.07%  '
$ws.Range("D16").Value = '87.20'
$ws.Range("E16").Value = '  -5.18%  '
$ws.Range("D17").Value = '1.008'
$ws.Range("E17").Value = '  +0.57%  '
$ws.Range("D18").Value = '0.000008713'
$ws.Range("D19").Value = '1.006'
$ws.Range("E19").Value = '  +0.56%  '
$ws.Range("D20").Value = '26.986.27'
$ws.Range("E20").Value = '  -1.27%  '
$ws.Range("D21").Value = '14.50'
$ws.Range("E21").Value = '  -1.96%  '
$ws.Range("D22").Value = '5.009'
$ws.Range("E22").Value = '  -2.13%  '
$ws.Range("D23").Value = '2.060.54'
$ws.Range("E23").Value = '  -0.85%  '
$ws.Range("E24").Value = '  -0.45%  '
$ws.Range("D25").Value = '2.015'
$ws.Range("E25").Value = '  +6.76%  '
$ws.Range("D26").Value = '151.45'
$ws.Range("E26").Value = '  -0.57%  '
$ws.Range("D27").Value = '18.21'
$ws.Range("E27").Value = '  -0.77%  '
$ws.Range("D28").Value = '1.964'
$ws.Range("E28").Value = '  -5.03%  '
$ws.Range("D29").Value = '113.98'
$ws.Range("E29").Value = '  -1.64%  '
$ws.Range("D30").Value = '4.932'
$ws.Range("E30").Value = '  -3.53%  '
$ws.Range("D31").Value = '0.08805'
$ws.Range("E31").Value = '  -0.38%  '
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").Value = '0.7510'
$ws.Range("E32").Value = '  -1.62%  '
$ws.Range("B33").Value = 'HuobiToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D33").Value = '2.986'
$ws.Range("E33").Value = '  -0.33%  '
$ws.Range("D34").Value = '4.487'
$ws.Range("E34").Value = '  -0.06%  '
$ws.Range("E35").Value = '  -2.78%  '
$ws.Range("D36").Value = '2.560'
$ws.Range("E36").Value = '  -1.84%  '
$ws.Range("D37").Value = '1.090'
$ws.Range("E37").Value = '  +1.18%  '
$ws.Range("D38").Value = '0.01949'
$ws.Range("E38").Value = '  -0.58%  '
$ws.Range("D39").Value = '0.05155'
$ws.Range("E39").Value = '  -1.10%  '
$ws.Range("D40").Value = '2.893'
$ws.Range("E40").Value = '  -3.23%  '
$ws.Range("D41").Value = '6.951'
$ws.Range("E41").Value = '  -1.13%  '
$ws.Range("D42").Value = '0.4991'
$ws.Range("E42").Value = '  -3.10%  '
$ws.Range("E43").Value = '  -2.25%  '
$ws.Range("D44").Value = '8.327'
$ws.Range("E44").Value = '  -0.19%  '
$ws.Range("D45").Value = '0.4683'
$ws.Range("E45").Value = '  -3.10%  '
$ws.Range("D46").Value = '1.006'
$ws.Range("E46").Value = '  +0.49%  '
$ws.Range("D47").Value = '10.12'
$ws.Range("E47").Value = '  -1.46%  '
$ws.Range("D48").Value = '102.20'
$ws.Range("E48").Value = '  -1.01%  '
$ws.Range("D49").Value = '1.612'
$ws.Range("E49").Value = '  -2.34%  '
$ws.Range("D50").Value = '0.06110'
$ws.Range("E50").Value = '  -1.83%  '
$ws.Range("D51").Value = '64.46'
$ws.Range("E51").Value = '  -1.67%  '

$dRange.ClearFormats()
